$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Sales column (B2:B25) with new values
$newValues = @(344, 444, 428, 612, 558, 301, 519, 457, 439, 547, 600, 572, 422, 476, 520, 610, 422, 465, 376, 383, 626, 642, 457, 546)

$row = 2
foreach ($val in $newValues) {
    $ws.Cells.Item($row, 2).Value = $val
    $row++
}

# Update the active selection to C9
$ws.Range("C9").Select()
